{"js": "// Update the five rows of division problems in the worksheet table.\n// The table has 20 rows total; every 4th row (0, 4, 8, 12, 16) holds 5\n// division expressions, the rows in between are blank answer rows.\n// We overwrite the populated rows' cell text in place (left-to-right,\n// top-to-bottom), which keeps each run's existing formatting intact.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newRows = {\n  0: [\"82\u00f78=\", \"52\u00f79=\", \"56\u00f76=\", \"62\u00f77=\", \"26\u00f76=\"],\n  4: [\"78\u00f73=\", \"74\u00f75=\", \"92\u00f76=\", \"36\u00f79=\", \"94\u00f78=\"],\n  8: [\"52\u00f79=\", \"74\u00f79=\", \"36\u00f75=\", \"50\u00f74=\", \"18\u00f77=\"],\n  12: [\"83\u00f76=\", \"83\u00f77=\", \"84\u00f75=\", \"19\u00f72=\", \"92\u00f75=\"],\n  16: [\"48\u00f74=\", \"48\u00f77=\", \"88\u00f72=\", \"54\u00f77=\", \"25\u00f76=\"],\n};\n\nconst values = table.values;\nfor (const rowIndexStr of Object.keys(newRows)) {\n  const rowIndex = Number(rowIndexStr);\n  values[rowIndex] = newRows[rowIndex];\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the five rows of division problems in the worksheet table.\n# The table has 20 rows total; every 4th row (1, 5, 9, 13, 17 in 1-based\n# COM indexing) holds 5 division expressions, the rows in between are\n# blank answer rows. We overwrite each populated cell's text in place,\n# which keeps the existing run formatting (font/size) intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRows = @{\n    1  = @(\"82\u00f78=\", \"52\u00f79=\", \"56\u00f76=\", \"62\u00f77=\", \"26\u00f76=\")\n    5  = @(\"78\u00f73=\", \"74\u00f75=\", \"92\u00f76=\", \"36\u00f79=\", \"94\u00f78=\")\n    9  = @(\"52\u00f79=\", \"74\u00f79=\", \"36\u00f75=\", \"50\u00f74=\", \"18\u00f77=\")\n    13 = @(\"83\u00f76=\", \"83\u00f77=\", \"84\u00f75=\", \"19\u00f72=\", \"92\u00f75=\")\n    17 = @(\"48\u00f74=\", \"48\u00f77=\", \"88\u00f72=\", \"54\u00f77=\", \"25\u00f76=\")\n}\n\nforeach ($rowIndex in $newRows.Keys) {\n    $values = $newRows[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
